# Upload the actual water-taste-test data (replacing the earlier placeholder
# values) and add a new "X" sample row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "name" column for samples A-E with the real water brands.
# Written in this order so new shared-strings get appended in the same
# sequence as the authoritative workbook.
$ws.Range("B6").Value = "Well Water"
$ws.Range("B2").Value = "Evian"
$ws.Range("B5").Value = "Smart Water"
$ws.Range("B4").Value = "Press Building Water"
$ws.Range("B3").Value = "Fiji"

# Add the new sample "X", also sourced from Well Water.
$ws.Range("A7").Value = "X"
$ws.Range("B7").Value = "Well Water"

# Leave the selection where data entry finished (one row below the last
# entered row), matching normal Excel behavior after typing + Enter.
$ws.Range("B8").Select() | Out-Null
